# Refresh the cryptocurrency price / 1h-volume table with the latest
# scraped values (Tue Dec 19 05:55:05 UTC 2023, GitHub Actions run).
#
# Numeric-looking price strings (e.g. "76.21") are written through
# Range.NumberFormat = "@" so Excel keeps them as literal text instead
# of auto-converting to numbers, matching the source data which stores
# every Price/Volume cell as text. Style is reset back to "Normal"
# afterwards so no stray number-format style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.928.07"
$ws.Range("E2").Value = "  +4.25%  "
$ws.Range("D3").Value = "2.241.72"
$ws.Range("E3").Value = "  +3.36%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.21"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.72%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +6.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("D15").Value = "2.557.26"
$ws.Range("E15").Value = "  +2.53%  "
$ws.Range("E16").Value = "  +5.43%  "
$ws.Range("D17").Value = "2.243.81"
$ws.Range("E17").Value = "  +3.88%  "
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "42.821.79"
$ws.Range("E19").Value = "  +4.44%  "
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.87%  "
$ws.Range("E24").Value = "  +14.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "231.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -4.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +23.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.11%  "
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0794"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.66%  "
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.76%  "
$ws.Range("E39").Value = "  +15.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.94%  "
$ws.Range("E41").Value = "  +3.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("E43").Value = "  +6.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0989"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.441"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.96%  "
$ws.Range("E51").Value = "  +1.41%  "

# Rows 44-45: the ranking flipped Aave <-> MultiversX, each with
# refreshed price/volume figures.
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.76%  "

$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "59.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.53%  "
